$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44211
$ws.Cells.Item(2, 9).Value = 'Segunda'
$ws.Cells.Item(2, 11).Value = 4500
$ws.Cells.Item(2, 12).Value = 5000
$ws.Cells.Item(2, 13).Value = 4750
$ws.Cells.Item(2, 16).Value = 475

$ws.Cells.Item(3, 4).Value = 44377
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 17000
$ws.Cells.Item(3, 12).Value = 18000
$ws.Cells.Item(3, 13).Value = 17600
$ws.Cells.Item(3, 16).Value = 978

$ws.Cells.Item(4, 4).Value = 44405
$ws.Cells.Item(4, 9).Value = 'Segunda'
$ws.Cells.Item(4, 10).Value = 140
$ws.Cells.Item(4, 13).Value = 17500
$ws.Cells.Item(4, 16).Value = 972

$ws.Cells.Item(5, 4).Value = 44363
$ws.Cells.Item(5, 10).Value = 140
$ws.Cells.Item(5, 11).Value = 14000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 14500
$ws.Cells.Item(5, 16).Value = 806

$ws.Cells.Item(6, 4).Value = 44221
$ws.Cells.Item(6, 8).Value = 'Cultivar XV región'
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 140
$ws.Cells.Item(6, 11).Value = 5000
$ws.Cells.Item(6, 12).Value = 6000
$ws.Cells.Item(6, 13).Value = 5500
$ws.Cells.Item(6, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(6, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(6, 16).Value = 550
$ws.Cells.Item(6, 17).Value = 10

$ws.Cells.Item(7, 4).Value = 44433
$ws.Cells.Item(7, 9).Value = 'Segunda'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17500
$ws.Cells.Item(7, 16).Value = 972

$ws.Cells.Item(8, 4).Value = 44433
$ws.Cells.Item(8, 9).Value = 'Tercera'
$ws.Cells.Item(8, 10).Value = 120
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 14500
$ws.Cells.Item(8, 16).Value = 806

$ws.Cells.Item(9, 4).Value = 44398
$ws.Cells.Item(9, 9).Value = 'Primera'
$ws.Cells.Item(9, 11).Value = 17000
$ws.Cells.Item(9, 12).Value = 18000
$ws.Cells.Item(9, 13).Value = 17500
$ws.Cells.Item(9, 16).Value = 972

$ws.Cells.Item(10, 4).Value = 44398
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 16000
$ws.Cells.Item(10, 13).Value = 15500
$ws.Cells.Item(10, 16).Value = 861

$ws.Cells.Item(11, 4).Value = 44391
$ws.Cells.Item(11, 9).Value = 'Segunda'
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 16000
$ws.Cells.Item(11, 13).Value = 15500
$ws.Cells.Item(11, 16).Value = 861

$ws.Cells.Item(12, 4).Value = 44435
$ws.Cells.Item(12, 8).Value = 'Cultivar IV Región'
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 17000
$ws.Cells.Item(12, 12).Value = 18000
$ws.Cells.Item(12, 13).Value = 17500
$ws.Cells.Item(12, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(12, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(12, 16).Value = 972
$ws.Cells.Item(12, 17).Value = 18

$ws.Cells.Item(13, 4).Value = 44435
$ws.Cells.Item(13, 9).Value = 'Tercera'
$ws.Cells.Item(13, 10).Value = 120
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14500
$ws.Cells.Item(13, 16).Value = 806

$ws.Cells.Item(14, 4).Value = 44412
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 150
$ws.Cells.Item(14, 11).Value = 17000
$ws.Cells.Item(14, 12).Value = 18000
$ws.Cells.Item(14, 13).Value = 17500
$ws.Cells.Item(14, 16).Value = 972
